# feat(commands): add the command cutsendmail to cut a file in multiple
# tabs which are eventually sent by mail to their owner.
#
# The worksheet "sheet1" repeats a 4-column block (Alain/Henri/Tony/
# Dulcinee header in row 1, OUI/NON answers in rows 2-9) many times across
# the row, followed by a trailing "email" column and an empty column.
# This change extends the repeated block by 5 more repetitions (20 columns)
# so the per-recipient tabs created by the new "cutsendmail" command line
# up on 4-column boundaries. The trailing email / empty columns are pushed
# further right by the same 20 columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$insertCount = 20
$blockWidth = 4

# Column immediately to the right of the last existing repeated block
# (right before the trailing "email" / empty columns) -- this is where the
# new columns need to be inserted.
$insertFirstCol = "JA"
$insertLastCol = "JT"

# The last existing 4-column block, which we replicate into the newly
# inserted columns.
$sourceBlock = "IW1:IZ9"

# Shift everything from JA: onward (email + empty columns, and the used
# dimension) to the right by inserting 20 blank columns.
$insertRange = $ws.Range($insertFirstCol + "1:" + $insertLastCol + "1")
$insertRange.EntireColumn.Insert()

# Fill the freshly inserted columns by tiling the last 4-column block five
# times (20 columns), carrying over both values and formatting.
$ws.Range($sourceBlock).Copy()

$destStarts = @("JA1", "JE1", "JI1", "JM1", "JQ1")
foreach ($startCell in $destStarts) {
    $ws.Range($startCell).PasteSpecial(-4104)
}

$excel.CutCopyMode = 0
